$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "Electrode Locations" column (C) derived from the file name in
# column A (e.g. "A11_monopolar_10V_1kHz.txt" -> "A11"), then sort the whole
# data block (rows 2-70) by that electrode location using a natural sort
# (letter prefix, then numeric suffix) so the order runs A1 .. O15 instead
# of the old plain alphabetical file-name order.
# ---------------------------------------------------------------------------

$lastRow = 70

# Header for the new column, matching the style already used by A1/B1
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Read the existing data (file name + value) into memory
$fileNames = @()
$values = @()
$locations = @()
$sortKeys = @()

$n = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $fileName = [string]$ws.Cells.Item($r, 1).Value()
    $value = $ws.Cells.Item($r, 2).Value()

    # Electrode location is the token before the first underscore
    $loc = $fileName.Split("_")[0]

    # Split the location into its alphabetic prefix and numeric suffix so the
    # rows can be ordered naturally (A2, A3, A8 ... A11, A14) rather than
    # lexicographically (A11, A14, A2 ...). Zero-pad the number so a plain
    # string sort on "letters|number|index" yields the natural order, and
    # keep the original index as a tie-breaker/back-reference.
    $letterPart = ($loc -replace '[0-9]+$', '')
    $digitPart = [int]($loc -replace '^[A-Za-z]+', '')
    $paddedNum = "{0:D4}" -f $digitPart
    $paddedIdx = "{0:D4}" -f $n

    $fileNames += , $fileName
    $values += , $value
    $locations += , $loc
    $sortKeys += , "$letterPart|$paddedNum|$paddedIdx"

    $n++
}

# Natural sort via the composite string key
$sortedKeys = $sortKeys | Sort-Object

# Write the sorted data (and the derived location) back to the sheet. Rows
# that already sit at the correct destination keep their original A/B cells
# untouched (only the new column C is populated) so cells that don't need to
# move are not needlessly rewritten; rows that do move are written in full.
$r = 2
foreach ($key in $sortedKeys) {
    $parts = $key.Split("|")
    $i = [int]$parts[2]
    $srcRow = $i + 2

    if ($srcRow -ne $r) {
        $ws.Cells.Item($r, 1).Value = $fileNames[$i]
        $ws.Cells.Item($r, 2).Value = $values[$i]
    }
    $ws.Cells.Item($r, 3).Value = $locations[$i]
    $r++
}
